$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to Text so numeric-looking price/percent strings are not
# auto-converted to numbers by Excel, then drop back to the default style so
# no stray number-format is left attached to the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.137.88'
$ws.Range("E2").Value = '  +2.36%  '

$ws.Range("D3").Value = '1.652.04'
$ws.Range("E3").Value = '  +2.16%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '213.52'
$ws.Range("E5").Value = '  +1.23%  '

$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("D8").Value = '23.56'
$ws.Range("E8").Value = '  +3.47%  '

$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("D11").Value = '0.0873'
$ws.Range("E11").Value = '  -1.39%  '

$ws.Range("D12").Value = '1.883.10'

$ws.Range("D13").Value = '1.654.15'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("D15").Value = '0.567'
$ws.Range("E15").Value = '  +3.22%  '

$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("D17").Value = '28.108.26'
$ws.Range("E17").Value = '  +2.30%  '

$ws.Range("D18").Value = '233.38'
$ws.Range("E18").Value = '  +0.96%  '

$ws.Range("D19").Value = '7.71'
$ws.Range("E19").Value = '  +2.51%  '

$ws.Range("D20").Value = '0.0₃0724'
$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("D21").Value = '0.998'
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("E22").Value = '  +5.37%  '

$ws.Range("D23").Value = '4.41'
$ws.Range("E23").Value = '  +2.88%  '

$ws.Range("E24").Value = '  +4.06%  '

$ws.Range("D25").Value = '152.06'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("D26").Value = '6.93'
$ws.Range("E26").Value = '  +1.29%  '

$ws.Range("D27").Value = '15.79'
$ws.Range("E27").Value = '  +1.72%  '

$ws.Range("E28").Value = '  +0.40%  '

$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.30%  '

$ws.Range("D30").Value = '1.19'
$ws.Range("E30").Value = '  +1.53%  '

$ws.Range("E31").Value = '  +0.33%  '

$ws.Range("E32").Value = '  +2.60%  '

$ws.Range("D33").Value = '1.448.22'
$ws.Range("E33").Value = '  -1.37%  '

$ws.Range("E34").Value = '  +0.33%  '

$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +2.53%  '

$ws.Range("E36").Value = '  -0.83%  '

$ws.Range("D37").Value = '0.892'
$ws.Range("E37").Value = '  +3.72%  '

$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Value = '0.923'
$ws.Range("E40").Value = '  -2.07%  '

$ws.Range("D41").Value = '69.51'
$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("E42").Value = '  +3.47%  '

$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  -0.38%  '

$ws.Range("D44").Value = '2.46'
$ws.Range("E44").Value = '  -0.62%  '

$ws.Range("D45").Value = '1.82'
$ws.Range("E45").Value = '  +5.93%  '

$ws.Range("B46").Value = 'MXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  +0.96%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '5.40'
$ws.Range("E47").Value = '  +2.64%  '

$ws.Range("D48").Value = '1.794.90'
$ws.Range("E48").Value = '  +1.98%  '

$ws.Range("D49").Value = '88.90'
$ws.Range("E49").Value = '  +2.64%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.101'
$ws.Range("E50").Value = '  +0.48%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.73'
$ws.Range("E51").Value = '  +0.81%  '

$ws.Range("D2:E51").Style = "Normal"
